# Exclude Wuhan row (CHN1, "Wuhan, China") from the derived data map sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Row 12 contains the Wuhan, China entry (CHN1 / ageband). Delete the entire
# row, shifting subsequent rows up.
$ws.Rows.Item(12).Delete()

# Reflect the post-edit selection, matching where the cursor ended up.
$ws.Range("A12:XFD12").Select()
